$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trial block for 2025-08-18 (row 13), plus its tally row (row 14)
$ws.Range("A13").Value = 20250818
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 5
$ws.Range("G13").Value = 6

$ws.Range("C14").Value = 17
$ws.Range("E14").Value = 33

# New trial block for 2025-08-20 (row 15), plus its tally row (row 16)
$ws.Range("A15").Value = 20250820
$ws.Range("B15").Value = 3
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = 5
$ws.Range("G15").Value = 6

$ws.Range("F16").Value = 24

# Match the author's final selection/active cell
$ws.Range("G16").Select()
